# Replace the "{Contact_Information}" placeholder with "{RSVP}" in the
# "Please RSVP by ..." paragraph.
#
# In the source XML this placeholder is spread across 5 separate runs:
#   "{"  / "Contact" / "_" / "Information" / "}"
# and the target state collapses the middle three runs ("Contact", "_",
# "Information") into a single run containing "RSVP", while the "{" and
# "}" runs (and the following sentence run) stay untouched.
#
# A plain Range.Text assignment in this host coalesces every run from the
# edited run through the end of the paragraph, which would also swallow
# the trailing "}" run and the sentence that follows it - more than the
# diff calls for. Range.InsertXML gives precise control over the OOXML
# that lands in the document, but it only behaves correctly (and safely)
# when applied to a Range spanning the *entire* paragraph (including its
# paragraph mark); using it on a sub-paragraph Range drops the rest of
# the paragraph's content. So: locate the whole paragraph, then replace
# it with an equivalent paragraph whose runs already reflect the desired
# split.

$d = $word.ActiveDocument

foreach ($para in $d.Paragraphs) {
    $full = $para.Range.Text
    if ($full.Contains("{Contact_Information}")) {

        $pPrXml = "<w:pPr><w:pStyle w:val=`"BodyText`"/></w:pPr>"

        $runsXml =
            "<w:r><w:t xml:space=`"preserve`">Please RSVP by </w:t></w:r>" +
            "<w:r w:rsidR=`"000E737A`"><w:t>{</w:t></w:r>" +
            "<w:r><w:t>RSVP</w:t></w:r>" +
            "<w:r w:rsidR=`"000E737A`"><w:t>_</w:t></w:r>" +
            "<w:r><w:t>Date</w:t></w:r>" +
            "<w:r w:rsidR=`"000E737A`"><w:t>}</w:t></w:r>" +
            "<w:r><w:t xml:space=`"preserve`"> at </w:t></w:r>" +
            "<w:r w:rsidR=`"000E737A`"><w:t>{</w:t></w:r>" +
            "<w:r><w:t>RSVP</w:t></w:r>" +
            "<w:r w:rsidR=`"000E737A`"><w:t>}</w:t></w:r>" +
            "<w:r><w:t>. Your presence will make this event even more special.</w:t></w:r>"

        $paraXml = "<w:p xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" " +
                   "xmlns:w14=`"http://schemas.microsoft.com/office/word/2010/wordml`" " +
                   "w14:paraId=`"72D06F16`" w14:textId=`"64D7D71D`" w:rsidR=`"00A404AA`" w:rsidRDefault=`"00000000`">" +
                   $pPrXml + $runsXml + "</w:p>"

        $para.Range.InsertXML($paraXml)
        break
    }
}
